$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 26 de Junio de 2020 a las 14:34'

$ws.Range("B4").Value = 2505196
$ws.Range("C4").Value = 608
$ws.Range("D4").Value = 1052392
$ws.Range("E4").Value = 1326006
$ws.Range("G4").Value = 18
$ws.Range("H4").Value = 126798

$ws.Range("B26").Value = 65137
$ws.Range("C26").Value = 303
$ws.Range("G26").Value = 50
$ws.Range("H26").Value = 5280

$ws.Range("B29").Value = 60713
$ws.Range("C29").Value = 331
$ws.Range("D29").Value = 42689
$ws.Range("E29").Value = 17651
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 373

$ws.Range("B33").Value = 50005
$ws.Range("C33").Value = 91
$ws.Range("G33").Value = 3
$ws.Range("H33").Value = 6103

$ws.Range("A35").Value = 'Kuwait'
$ws.Range("B35").Value = 43703
$ws.Range("C35").Value = 915
$ws.Range("D35").Value = 33969
$ws.Range("E35").Value = 9393
$ws.Range("G35").Value = 2
$ws.Range("H35").Value = 341

$ws.Range("A36").Value = 'Singapur'
$ws.Range("B36").Value = 42955
$ws.Range("C36").Value = 219
$ws.Range("D36").Value = 36604
$ws.Range("E36").Value = 6325
$ws.Range("H36").Value = 26

$ws.Range("B63").Value = 12675
$ws.Range("C63").Value = 39
$ws.Range("D63").Value = 11508
$ws.Range("E63").Value = 563
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 604

$ws.Range("B87").Value = 4643
$ws.Range("C87").Value = 8
$ws.Range("D87").Value = 4348
$ws.Range("E87").Value = 243

$ws.Range("B93").Value = 3935
$ws.Range("C93").Value = 139
$ws.Range("D93").Value = 2338
$ws.Range("E93").Value = 1419
$ws.Range("G93").Value = 3
$ws.Range("H93").Value = 178

$ws.Range("A101").Value = 'Croacia'
$ws.Range("B101").Value = 2539
$ws.Range("C101").Value = 56
$ws.Range("D101").Value = 2150
$ws.Range("E101").Value = 282
$ws.Range("H101").Value = 107

$ws.Range("A102").Value = 'Mayotte'
$ws.Range("B102").Value = 2508
$ws.Range("D102").Value = 2218
$ws.Range("E102").Value = 258
$ws.Range("H102").Value = 32

$ws.Range("D108").Value = 1619
$ws.Range("E108").Value = 380

$ws.Range("A112").Value = 'Madagascar'
$ws.Range("B112").Value = 1922
$ws.Range("C112").Value = 93
$ws.Range("D112").Value = 862
$ws.Range("E112").Value = 1044
$ws.Range("H112").Value = 16

$ws.Range("A113").Value = 'Islandia'
$ws.Range("B113").Value = 1830
$ws.Range("D113").Value = 1811
$ws.Range("E113").Value = 9
$ws.Range("H113").Value = 10

$ws.Range("B131").Value = 1053
$ws.Range("C131").Value = 36
$ws.Range("D131").Value = 292
$ws.Range("E131").Value = 747

$ws.Range("D147").Value = 656
$ws.Range("E147").Value = 0

$ws.Range("A202").Value = 'Dominica'

$ws.Range("A203").Value = 'Fiyi'

$ws.Range("A208").Value = 'Islas Malvinas'

$ws.Range("A209").Value = 'Groenlandia'

$ws.Range("A211").Value = 'Montserrat'
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

$ws.Range("A212").Value = 'Seychelles'
$ws.Range("D212").Value = 11
$ws.Range("H212").Value = 0
